$wb = $excel.ActiveWorkbook

# --- Sheet "DBD" (sheet1): CreateDate/LastUpdate column type DATE -> TIMESTAMP,
#     and clear the now-meaningless length value next to it.
$ws1 = $wb.Worksheets.Item("DBD")

$ws1.Range("D12").Value = "TIMESTAMP"
$ws1.Range("E12").ClearContents()

$ws1.Range("D14").Value = "TIMESTAMP"
$ws1.Range("E14").ClearContents()

# Leave the cursor where the author last left it before saving.
[void]$ws1.Range("E14").Select()

# --- Sheet "DBS" (sheet2): remove the stray empty formatted row 3
$ws2 = $wb.Worksheets.Item("DBS")
$ws2.Rows.Item(3).Delete()
